$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.294.00'
$ws.Range('E2').Value = '  -1.95%  '

$ws.Range('D3').Value = '2.526.94'
$ws.Range('E3').Value = '  -1.78%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.76'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.04%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.565'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.57%  '

$ws.Range('E8').Value = '  +0.10%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.53%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.44'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.85%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0798'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.77%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.30'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.24%  '

$ws.Range('E13').Value = '  -0.03%  '

$ws.Range('D14').Value = '2.920.76'
$ws.Range('E14').Value = '  -1.60%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.80'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.51%  '

$ws.Range('D16').Value = '2.461.39'
$ws.Range('E16').Value = '  -2.53%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.830'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.13%  '

$ws.Range('D18').Value = '42.309.16'
$ws.Range('E18').Value = '  -2.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.75'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.99%  '

$ws.Range('D20').Value = '0.0₃0944'
$ws.Range('E20').Value = '  -2.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.13'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.02%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.71'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.42%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.12'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.64%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.87'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.65%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.03'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.70%  '

$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.03'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.98%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.34'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.70%  '

$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.20'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.49%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.09'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.42'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.54%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.82'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +16.12%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.66'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.91%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0791'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.28%  '

$ws.Range('E35').Value = '  -3.07%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.02'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -5.49%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.07'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.98%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.14'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -7.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.110'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.95%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.117'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.83%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.23'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +9.18%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.37'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.67%  '

$ws.Range('E43').Value = '  +0.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.26'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.43%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0295'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.10%  '

$ws.Range('D46').Value = '1.953.05'
$ws.Range('E46').Value = '  -2.21%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.87'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.16%  '

$ws.Range('D48').Value = '2.777.45'
$ws.Range('E48').Value = '  -1.55%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '80.39'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.69%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.190'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.70%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.82'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.24%  '
